# "updated minutes of meeting"
# Appends the 27/09/2023 meeting entry (row 11) to the meeting diary table
# on Sheet1, matching the formatting already used by the rows above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new meeting-diary row -------------------------------------------------
# Columns: A=Date, B=Time start, C=Time end, D=Members present, E=Discussions
# Set the values first (while the cell still has no special number format)
# so the text in E11 - which starts with "-" - keeps the same "treat as
# text" quoting behaviour as the cells above it once we copy their format.
$ws.Range("A11").Value = 45196
$ws.Range("B11").Value = 0.89583333333333337
$ws.Range("C11").Value = 0.94444444444444453
$ws.Range("D11").Value = "All"
$ws.Range("E11").Value = "- Data quality check finalised`n- Data cleansing to be for Q3 part 2`n- Deadline for Q3 part 2 is 29/09/2023`n- Next meeting on 29/09/2023"

# Copy the formatting (date/time number formats, wrap text, etc.) from the
# row above so row 11 looks like the rest of the table, without disturbing
# the values we just entered.
$ws.Range("A10:E10").Copy()
$ws.Range("A11:E11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Resize the row to fit the (longer) wrapped discussion text.
$ws.Rows.Item(11).RowHeight = 62.4

# Leave the selection where the author's last save left it.
[void]$ws.Range("E13").Select()
